$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column N currently duplicates "TOTAL EMPLEADO" (same as G2); update it
# to "TOTAL PATRON" so the PATRON section (H1:N1 merged header) has its own
# total label, adding a new shared string in the process.
$ws.Range("N2").Value = "TOTAL PATRON"

# Update the active cell / selection to match the saved view state.
$ws.Range("K12").Select()
